# Update expected outputs using 10-13-22 PEARS schema
$wb = $excel.ActiveWorkbook

# --- Sheet: Unique Sites and Reach by Goal ---
$ws1 = $wb.Worksheets.Item("Unique Sites and Reach by Goal")
$ws1.Range("D2").Value = 7711746
$ws1.Range("D3").Value = 32881343
$ws1.Range("D4").Value = 6986656
$ws1.Range("D5").Value = 11370313
$ws1.Range("D6").Value = 42022897
$ws1.Range("D7").Value = 9941231
$ws1.Range("D8").Value = 13053623
$ws1.Range("D9").Value = 55462894
$ws1.Range("D10").Value = 10756157
$ws1.Range("C11").Value = 1193
$ws1.Range("D11").Value = 18382612
$ws1.Range("C12").Value = 484
$ws1.Range("D12").Value = 76978041
$ws1.Range("C13").Value = 617
$ws1.Range("D13").Value = 13049556

# --- Sheet: Direct Education Demographics ---
$ws2 = $wb.Worksheets.Item("Direct Education Demographics")
$ws2.Range("D4").Value = 0
$ws2.Range("C5").Value = 32
$ws2.Range("C6").Value = 14196
$ws2.Range("D6").Value = 1
$ws2.Range("D7").Value = 0
$ws2.Range("D8").Value = 1
$ws2.Range("C9").Value = 75
$ws2.Range("C11").Value = 6762
$ws2.Range("D11").Value = 0
$ws2.Range("C12").Value = 38
$ws2.Range("C13").Value = 19597
$ws2.Range("D13").Value = 1
$ws2.Range("D14").Value = 0
$ws2.Range("D15").Value = 1
$ws2.Range("C16").Value = 107
$ws2.Range("C18").Value = 8984
$ws2.Range("D18").Value = 0
$ws2.Range("C19").Value = 57
$ws2.Range("C20").Value = 25033
$ws2.Range("D20").Value = 1
$ws2.Range("D21").Value = 0
$ws2.Range("D22").Value = 1
$ws2.Range("C23").Value = 147
$ws2.Range("C25").Value = 14207
$ws2.Range("D25").Value = 0
$ws2.Range("C26").Value = 78
$ws2.Range("C27").Value = 33323
$ws2.Range("D27").Value = 1
$ws2.Range("D28").Value = 0
$ws2.Range("C29").Value = 40665
$ws2.Range("D29").Value = 1

# --- Sheet: RE-AIM Reach ---
$ws3 = $wb.Worksheets.Item("RE-AIM Reach")
$ws3.Range("B2").Value = 3220856
$ws3.Range("C2").Value = 5913807
$ws3.Range("E2").Value = 4749021107
$ws3.Range("F2").Value = 44461815
$ws3.Range("B3").Value = 4680171
$ws3.Range("C3").Value = 12865665
$ws3.Range("E3").Value = 8202568897
$ws3.Range("F3").Value = 59035828
$ws3.Range("B4").Value = 5946691
$ws3.Range("C4").Value = 18173965
$ws3.Range("D4").Value = 12216
$ws3.Range("E4").Value = 10063526541
$ws3.Range("F4").Value = 73781154
$ws3.Range("B5").Value = 8152156
$ws3.Range("C5").Value = 21191831
$ws3.Range("D5").Value = 14224
$ws3.Range("E5").Value = 11899487801
$ws3.Range("F5").Value = 101174712

# --- Sheet: RE-AIM Adoption ---
$ws4 = $wb.Worksheets.Item("RE-AIM Adoption")
$ws4.Range("B5").Value = 1340
